$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'320.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-3.24%"
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'42.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-6.21%"
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = "'5.203"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-6.66%"
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'0.08172"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-2.26%"
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'4.310"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-3.15%"
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'1.811"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-13.71%"
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'0.9335"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-5.66%"
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'0.1108"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-7.55%"
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'0.1860"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-3.64%"
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'0.04714"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'0.94%"
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'0.09377"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-5.42%"
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'7.431"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-27.92%"
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'0.1058"
$ws.Range("D14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'0.001289"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.36%"
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'0.005804"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.33%"
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = "'3.356"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-1.15%"
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").Value = "'2.508"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-1.49%"
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").Value = "'0.3348"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.31%"
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("B20").Value = 'ProBitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D20").Value = "'0.1381"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'0.63%"
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("B21").Value = 'ZBToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D21").Value = "'0.2548"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-0.64%"
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("B22").Value = 'CoinExToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D22").Value = "'0.04147"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-0.14%"
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("B23").Value = 'BitKan'
$ws.Range("C23").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D23").Value = "'0.001247"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-3.52%"
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("B24").Value = 'HotbitToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D24").Value = "'0.004299"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-5.17%"
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("E25").Value = "'-7.69%"
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'0.0002982"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-20.37%"
$ws.Range("E26").Style = "Normal"

# Row 38
$ws.Range("D38").Value = "'0.02693"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'-0.36%"
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'0.05545"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-3.64%"
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'0.008128"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'2.83%"
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'0.1399"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-2.46%"
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'0.006545"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-14.08%"
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("E43").Value = "'3.23%"
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'0.008254"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-7.56%"
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'0.3492"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'2.50%"
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'0.00006928"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-1.70%"
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("E47").Value = "'0.00%"
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'0.003347"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-2.25%"
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "'0.003534"
$ws.Range("D49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.00%"
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.00%"
$ws.Range("E51").Style = "Normal"
